# Redo gamemode config / add KitSync + ConfigSync net calls to the "netcalls" table.
#
# Inserts three new rows (113-115) into Table2 on the "netcalls" sheet, pushing
# the existing SendReportInvocation / ReceiveInvocationResponse rows down to
# 116-117, then expands the table to cover the new range and updates the
# selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new rows right before the old row 113.
$ws.Rows("113:115").Insert() | Out-Null

# --- Row 113: MulticastKitUpdated (KitSync.NetCalls, FROM_EITHER, Kit kit) ---
$ws.Range("C113").Value = "KitSync.NetCalls"
$ws.Range("E113").Value = "Kit kit"
$ws.Range("A113").Value = "MulticastKitUpdated"
$ws.Range("D113").Value = "FROM_EITHER"
$ws.Range("B113").Value = 3001

# --- Row 114: SendSingleProperty (ConfigSync.NetCalls, FROM_SERVER, PropertyValue property) ---
$ws.Range("C114").Value = "ConfigSync.NetCalls"
$ws.Range("D114").Value = "FROM_SERVER"
$ws.Range("B114").Value = 3002

# --- Row 115: ReceiveSyncPacket (ConfigSync.NetCalls, FROM_SERVER, SyncPacket packet) ---
$ws.Range("E115").Value = "SyncPacket packet"
$ws.Range("A115").Value = "ReceiveSyncPacket"
$ws.Range("C115").Value = "ConfigSync.NetCalls"
$ws.Range("D115").Value = "FROM_SERVER"
$ws.Range("B115").Value = 3003

# Fill in the two remaining cells of row 114 last.
$ws.Range("A114").Value = "SendSingleProperty"
$ws.Range("E114").Value = "PropertyValue property"

# Expand Table2 ("NetCall Field Name"/"ID"/"Containing Type"/"Direction"/"Method Signature")
# to include the three new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E117")) | Out-Null

# Match the author's final selection.
$ws.Range("A115").Select() | Out-Null
